$d = $word.ActiveDocument

# --- Locate the insertion point: right after "research position ",
#     right before "available within your group..." ---
$finder = $d.Content
$finder.Find.Execute("research position ")
$splitPoint = $finder.End

# End of the paragraph ("I am applying ...") so we only touch this
# paragraph's own run and not the following differently-formatted run.
$paraEnd = $d.Paragraphs(2).Range.End

# --- Cleanly split the existing run into a "prefix" run ending at
#     $splitPoint and a "tail" run, without disturbing neighbouring
#     runs that already have different formatting/rsid. Toggling a
#     character property on and back off forces the run boundary
#     without altering visible formatting. ---
$tail = $d.Range($splitPoint, $paraEnd)
$tail.Font.Bold = 1
$tail.Font.Bold = 0

# --- Move the (singleton) "_GoBack" bookmark to the split point; this
#     both inserts bookmarkStart/bookmarkEnd here and removes it from
#     its previous location automatically, exactly as real Word does
#     when tracking the most recent edit location. ---
$mark = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $mark)

# --- Insert the new "(PV1732) " text right before the split point.
#     Temporarily bolding the prefix run keeps the new text from being
#     folded back into it; we un-bold both afterwards. ---
$prefix = $d.Range(16, $splitPoint)
$prefix.Font.Bold = 1

$newText = "(PV1732) "
$ins = $d.Range($splitPoint, $splitPoint)
$ins.InsertBefore($newText)

$newLen = $newText.Length
$prefix2 = $d.Range(16, $splitPoint)
$prefix2.Font.Bold = 0
$newRun = $d.Range($splitPoint, $splitPoint + $newLen)
$newRun.Font.Bold = 0
